$d = $word.ActiveDocument

# Update the date heading
$d.Content.Find.Execute("2024-06-22 Saturday", $true, $false, $false, $false, $false, $true, 1, $false, "2024-06-23 Sunday", 2)

# Update the division problems in the table (5 data rows at table rows 1,5,9,13,17; 5 columns each)
$t = $d.Tables.Item(1)

$t.Cell(1, 1).Range.Text = "41÷4="
$t.Cell(1, 2).Range.Text = "86÷5="
$t.Cell(1, 3).Range.Text = "76÷6="
$t.Cell(1, 4).Range.Text = "50÷6="
$t.Cell(1, 5).Range.Text = "33÷4="

# $t.Cell(5, 1) stays "75÷7=" (unchanged)
$t.Cell(5, 2).Range.Text = "16÷6="
$t.Cell(5, 3).Range.Text = "69÷3="
$t.Cell(5, 4).Range.Text = "27÷7="
$t.Cell(5, 5).Range.Text = "30÷9="

$t.Cell(9, 1).Range.Text = "29÷6="
$t.Cell(9, 2).Range.Text = "39÷9="
$t.Cell(9, 3).Range.Text = "16÷8="
$t.Cell(9, 4).Range.Text = "43÷4="
$t.Cell(9, 5).Range.Text = "21÷6="

$t.Cell(13, 1).Range.Text = "20÷3="
$t.Cell(13, 2).Range.Text = "69÷8="
$t.Cell(13, 3).Range.Text = "61÷2="
$t.Cell(13, 4).Range.Text = "50÷3="
$t.Cell(13, 5).Range.Text = "17÷2="

$t.Cell(17, 1).Range.Text = "42÷8="
$t.Cell(17, 2).Range.Text = "56÷3="
$t.Cell(17, 3).Range.Text = "32÷4="
$t.Cell(17, 4).Range.Text = "65÷8="
$t.Cell(17, 5).Range.Text = "40÷8="
